$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ColorScheme
$cs.Colors(1).RGB  = 0x000000  # dk1
$cs.Colors(2).RGB  = 0xFFFFFF  # lt1
$cs.Colors(3).RGB  = 0x6A5444  # dk2 (BGR of 44546A)
$cs.Colors(4).RGB  = 0xE6E6E7  # lt2 (BGR of E7E6E6)
$cs.Colors(5).RGB  = 0xD59B5B  # accent1 BGR of 5B9BD5
$cs.Colors(6).RGB  = 0x317DED  # accent2 BGR of ED7D31
$cs.Colors(7).RGB  = 0xA5A5A5  # accent3
$cs.Colors(8).RGB  = 0x00C0FF  # accent4 BGR of FFC000
$cs.Colors(9).RGB  = 0xC47244  # accent5 BGR of 4472C4
$cs.Colors(10).RGB = 0x47AD70  # accent6 BGR of 70AD47
$cs.Colors(11).RGB = 0xC16305  # hlink BGR of 0563C1
$cs.Colors(12).RGB = 0x724F95  # folHlink BGR of 954F72
